$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.005433576865364477
$ws.Range("B3").Value = 0.7683170787553154
$ws.Range("B4").Value = 0.1004264215478384
$ws.Range("B5").Value = '[1.0, 0.01725367611630324, -0.07231940941169294, -0.05532376636399493, -0.09671975836559239, 0.13820335092151387, 0.06803179983886368, 0.05291314309029937, 0.038177140582319156, 0.045287223396612335, -0.002139298730657082, -0.033701805350904716, 0.011096970447720355, 0.037273268468067305, 0.07898957904792525, 0.03729759087333958, -0.015477061101819341, -0.009820360194695687, -0.012858705780243713, -0.011084155296903066]'

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1546161113014869
$ws.Range("B5").Value = 0.6297699631694609
$ws.Range("B6").Value = 0.0510167397912923
$ws.Range("B7").Value = 0.1539023900008555
$ws.Range("B8").Value = 0.6076610186296069
$ws.Range("B9").Value = 0.04893651470563919
$ws.Range("B10").Value = '[1.0, 0.12396570330166871, -0.005854818405963993, 0.018084495813276293, 0.0341108559206758, 0.02902000735434677, -0.003293867447770946, 0.032092444466914324, 0.02523830894547563, 0.021578400852515854, 0.020914166151904678, -0.0008959867308641242, -0.0012876800866132205, -0.019378398051696524, -0.028047799871448773, -0.026709684845678046, -0.03391741327343861, -0.015343108882375852, -0.03848653520976458, -0.007824847750583635]'

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.0001471735244537303
$ws.Range("B3").Value = 0.0203249893504086
$ws.Range("B4").Value = [double]"5.162627528144787E-08"
$ws.Range("B5").Value = 0.8212691572044912
$ws.Range("B6").Value = '[1.0, 0.013366203324305352, -0.06944102166179965, -0.08133344385288564, -0.11443852533211371, 0.1688209169813987, 0.06885018054144888, 0.05757761690568013, 0.04499659081015883, 0.037391599012683935, -0.032032469388477695, -0.04833174146775635, 0.0026685243080285417, 0.06346892333807547, 0.09504684112663411, 0.05386168862770744, -0.013553863190105713, -0.04363370454431479, -0.04460767103716553, -0.019516263290918164]'

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = 0.0002551511711441063
$ws.Range("B3").Value = 0.1126321227070231
$ws.Range("B4").Value = 0.01653271810613912
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = '[1.0, 0.006464966103609707, -0.06797729157429312, -0.07972534273645178, -0.11450882985376855, 0.16981153984242933, 0.06685367720488845, 0.05673594431810227, 0.044599939871101564, 0.037301345630700075, -0.031148517991422906, -0.04743482175499228, 0.0032607119732075453, 0.0625069138557058, 0.09368239812880454, 0.05354862240444152, -0.012670262002158742, -0.043065961766057086, -0.043496986170644285, -0.018942217487543597]'
$ws.Range("B7").Value = -0.01653271810157852

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.005507397691484082
$ws.Range("B3").Value = 0.1000333387772622
$ws.Range("B4").Value = 0.007569928075160842
$ws.Range("B5").Value = [double]"5.691837091256092E-10"
$ws.Range("B6").Value = '[1.0, 0.0134913442410164, -0.07170485565310351, -0.05450856765181642, -0.0967198119433175, 0.13869596726737968, 0.06743067916676847, 0.052596915790543965, 0.03789332005786615, 0.04539667342739283, -0.0018571992410704634, -0.033250816883565605, 0.011312951659183597, 0.03644726176443869, 0.07868018792776157, 0.03725634122715192, -0.015230342333643492, -0.009626361256997366, -0.01243516040698391, -0.01104977978542391]'
$ws.Range("B7").Value = -0.007569928123403419
$ws.Range("B9").Value = 0.7682907220672144
